# Update TPM-derived NATMI metrics for Col4a1-Itga1 sheet (rows 2-26).
# Values regenerated with new TPM normalization; only numeric measure columns change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 236.7946826666667
$ws.Range("H2").Value = 710.384048
$ws.Range("I2").Value = 0.7123899543147419
$ws.Range("J2").Value = 0.7240508783182559
$ws.Range("M2").Value = 32.21373866666666
$ws.Range("N2").Value = 96.641216
$ws.Range("O2").Value = 0.6812298485843321
$ws.Range("P2").Value = 0.7117693664123
$ws.Range("Q2").Value = 7628.042025080263
$ws.Range("R2").Value = 68652.37822572236
$ws.Range("S2").Value = 0.4853013007108309
$ws.Range("T2").Value = 0.5153572349108544

# Row 3
$ws.Range("G3").Value = 236.7946826666667
$ws.Range("H3").Value = 710.384048
$ws.Range("I3").Value = 0.7123899543147419
$ws.Range("J3").Value = 0.7240508783182559
$ws.Range("O3").Value = 0.1770389772624213
$ws.Range("P3").Value = 0.184975630381169
$ws.Range("Q3").Value = 1982.386358204032
$ws.Range("R3").Value = 17841.47722383629
$ws.Range("S3").Value = 0.1261207889239049
$ws.Range("T3").Value = 0.1339317676449585

# Row 4
$ws.Range("G4").Value = 236.7946826666667
$ws.Range("H4").Value = 710.384048
$ws.Range("I4").Value = 0.7123899543147419
$ws.Range("J4").Value = 0.7240508783182559
$ws.Range("M4").Value = 0.5484013333333334
$ws.Range("N4").Value = 1.645204
$ws.Range("O4").Value = 0.01159714372603029
$ws.Range("P4").Value = 0.01211704340205096
$ws.Range("Q4").Value = 129.8585197006436
$ws.Range("R4").Value = 1168.726677305792
$ws.Range("S4").Value = 0.008261688689168212
$ws.Range("T4").Value = 0.008773355917875425

# Row 5
$ws.Range("G5").Value = 236.7946826666667
$ws.Range("H5").Value = 710.384048
$ws.Range("I5").Value = 0.7123899543147419
$ws.Range("J5").Value = 0.7240508783182559
$ws.Range("M5").Value = 6.086836
$ws.Range("N5").Value = 12.173672
$ws.Range("O5").Value = 0.1287194389184112
$ws.Range("P5").Value = 0.08965995219214913
$ws.Range("Q5").Value = 1441.330399064043
$ws.Range("R5").Value = 8647.982394384257
$ws.Range("S5").Value = 0.09169843521050619
$ws.Range("T5").Value = 0.06491836713469841

# Row 6
$ws.Range("G6").Value = 236.7946826666667
$ws.Range("H6").Value = 710.384048
$ws.Range("I6").Value = 0.7123899543147419
$ws.Range("J6").Value = 0.7240508783182559
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.06689266666666667
$ws.Range("N6").Value = 0.200678
$ws.Range("O6").Value = 0.001414591508805173
$ws.Range("P6").Value = 0.001478007612330618
$ws.Range("Q6").Value = 15.83982777606044
$ws.Range("R6").Value = 142.558449984544
$ws.Range("S6").Value = 0.001007740780331739
$ws.Range("T6").Value = 0.001070152709869052

# Row 7
$ws.Range("I7").Value = 0.2358656137148928
$ws.Range("J7").Value = 0.2397264359793184
$ws.Range("M7").Value = 32.21373866666666
$ws.Range("N7").Value = 96.641216
$ws.Range("O7").Value = 0.6812298485843321
$ws.Range("P7").Value = 0.7117693664123
$ws.Range("Q7").Value = 2525.572971364003
$ws.Range("R7").Value = 22730.15674227603
$ws.Range("S7").Value = 0.160678696317247
$ws.Range("T7").Value = 0.1706299334492783

# Row 8
$ws.Range("I8").Value = 0.2358656137148928
$ws.Range("J8").Value = 0.2397264359793184
$ws.Range("O8").Value = 0.1770389772624213
$ws.Range("P8").Value = 0.184975630381169
$ws.Range("R8").Value = 5907.145306433112
$ws.Range("S8").Value = 0.04175740702345793
$ws.Range("T8").Value = 0.04434354861430538

# Row 9
$ws.Range("I9").Value = 0.2358656137148928
$ws.Range("J9").Value = 0.2397264359793184
$ws.Range("M9").Value = 0.5484013333333334
$ws.Range("N9").Value = 1.645204
$ws.Range("O9").Value = 0.01159714372603029
$ws.Range("P9").Value = 0.01211704340205096
$ws.Range("Q9").Value = 42.99493452958978
$ws.Range("R9").Value = 386.954410766308
$ws.Range("S9").Value = 0.002735367422279952
$ws.Range("T9").Value = 0.002904775629380392

# Row 10
$ws.Range("I10").Value = 0.2358656137148928
$ws.Range("J10").Value = 0.2397264359793184
$ws.Range("M10").Value = 6.086836
$ws.Range("N10").Value = 12.173672
$ws.Range("O10").Value = 0.1287194389184112
$ws.Range("P10").Value = 0.08965995219214913
$ws.Range("Q10").Value = 477.2109391522573
$ws.Range("R10").Value = 2863.265634913544
$ws.Range("S10").Value = 0.03036048945752771
$ws.Range("T10").Value = 0.02149386078909999

# Row 11
$ws.Range("I11").Value = 0.2358656137148928
$ws.Range("J11").Value = 0.2397264359793184
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.06689266666666667
$ws.Range("N11").Value = 0.200678
$ws.Range("O11").Value = 0.001414591508805173
$ws.Range("P11").Value = 0.001478007612330618
$ws.Range("Q11").Value = 5.244418000156222
$ws.Range("R11").Value = 47.19976200140599
$ws.Range("S11").Value = 0.0003336534943802082
$ws.Range("T11").Value = 0.0003543174972543212

# Row 12
$ws.Range("G12").Value = 0.6305213333333334
$ws.Range("H12").Value = 1.891564
$ws.Range("I12").Value = 0.001896905195629352
$ws.Range("J12").Value = 0.001927955138422806
$ws.Range("M12").Value = 32.21373866666666
$ws.Range("N12").Value = 96.641216
$ws.Range("O12").Value = 0.6812298485843321
$ws.Range("P12").Value = 0.7117693664123
$ws.Range("Q12").Value = 20.31144945575822
$ws.Range("R12").Value = 182.803045101824
$ws.Range("S12").Value = 0.001292228439197416
$ws.Range("T12").Value = 0.001372259407346539

# Row 13
$ws.Range("G13").Value = 0.6305213333333334
$ws.Range("H13").Value = 1.891564
$ws.Range("I13").Value = 0.001896905195629352
$ws.Range("J13").Value = 0.001927955138422806
$ws.Range("O13").Value = 0.1770389772624213
$ws.Range("P13").Value = 0.184975630381169
$ws.Range("Q13").Value = 5.278568233376001
$ws.Range("R13").Value = 47.50711410038401
$ws.Range("S13").Value = 0.0003358261557979935
$ws.Range("T13").Value = 0.0003566247170763726

# Row 14
$ws.Range("G14").Value = 0.6305213333333334
$ws.Range("H14").Value = 1.891564
$ws.Range("I14").Value = 0.001896905195629352
$ws.Range("J14").Value = 0.001927955138422806
$ws.Range("M14").Value = 0.5484013333333334
$ws.Range("N14").Value = 1.645204
$ws.Range("O14").Value = 0.01159714372603029
$ws.Range("P14").Value = 0.01211704340205096
$ws.Range("Q14").Value = 0.3457787398951112
$ws.Range("R14").Value = 3.112008659056
$ws.Range("S14").Value = 0.00002199868218836719
$ws.Range("T14").Value = 0.00002336111608947631

# Row 15
$ws.Range("G15").Value = 0.6305213333333334
$ws.Range("H15").Value = 1.891564
$ws.Range("I15").Value = 0.001896905195629352
$ws.Range("J15").Value = 0.001927955138422806
$ws.Range("M15").Value = 6.086836
$ws.Range("N15").Value = 12.173672
$ws.Range("O15").Value = 0.1287194389184112
$ws.Range("P15").Value = 0.08965995219214913
$ws.Range("Q15").Value = 3.837879950501334
$ws.Range("R15").Value = 23.027279703008
$ws.Range("S15").Value = 0.0002441685724628292
$ws.Range("T15").Value = 0.000172860365539597

# Row 16
$ws.Range("G16").Value = 0.6305213333333334
$ws.Range("H16").Value = 1.891564
$ws.Range("I16").Value = 0.001896905195629352
$ws.Range("J16").Value = 0.001927955138422806
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.06689266666666667
$ws.Range("N16").Value = 0.200678
$ws.Range("O16").Value = 0.001414591508805173
$ws.Range("P16").Value = 0.001478007612330618
$ws.Range("Q16").Value = 0.04217725337688889
$ws.Range("R16").Value = 0.379595280392
$ws.Range("S16").Value = 0.000002683345982745696
$ws.Range("T16").Value = 0.000002849532370820838

# Row 17
$ws.Range("G17").Value = 16.059769
$ws.Range("H17").Value = 32.119538
$ws.Range("I17").Value = 0.04831535056182164
$ws.Range("J17").Value = 0.032737474561192
$ws.Range("M17").Value = 32.21373866666666
$ws.Range("N17").Value = 96.641216
$ws.Range("O17").Value = 0.6812298485843321
$ws.Range("P17").Value = 0.7117693664123
$ws.Range("Q17").Value = 517.3452016130346
$ws.Range("R17").Value = 3104.071209678208
$ws.Range("S17").Value = 0.03291385894752868
$ws.Range("T17").Value = 0.02330153152635842

# Row 18
$ws.Range("G18").Value = 16.059769
$ws.Range("H18").Value = 32.119538
$ws.Range("I18").Value = 0.04831535056182164
$ws.Range("J18").Value = 0.032737474561192
$ws.Range("O18").Value = 0.1770389772624213
$ws.Range("P18").Value = 0.184975630381169
$ws.Range("Q18").Value = 134.448403245288
$ws.Range("R18").Value = 806.6904194717281
$ws.Range("S18").Value = 0.008553700249540254
$ws.Range("T18").Value = 0.006055634994043975

# Row 19
$ws.Range("G19").Value = 16.059769
$ws.Range("H19").Value = 32.119538
$ws.Range("I19").Value = 0.04831535056182164
$ws.Range("J19").Value = 0.032737474561192
$ws.Range("M19").Value = 0.5484013333333334
$ws.Range("N19").Value = 1.645204
$ws.Range("O19").Value = 0.01159714372603029
$ws.Range("P19").Value = 0.01211704340205096
$ws.Range("Q19").Value = 8.807198732625334
$ws.Range("R19").Value = 52.843192395752
$ws.Range("S19").Value = 0.0005603200646389837
$ws.Range("T19").Value = 0.0003966814001315027

# Row 20
$ws.Range("G20").Value = 16.059769
$ws.Range("H20").Value = 32.119538
$ws.Range("I20").Value = 0.04831535056182164
$ws.Range("J20").Value = 0.032737474561192
$ws.Range("M20").Value = 6.086836
$ws.Range("N20").Value = 12.173672
$ws.Range("O20").Value = 0.1287194389184112
$ws.Range("P20").Value = 0.08965995219214913
$ws.Range("Q20").Value = 97.75318010088399
$ws.Range("R20").Value = 391.012720403536
$ws.Range("S20").Value = 0.006219124815464026
$ws.Range("T20").Value = 0.002935240404048173

# Row 21
$ws.Range("G21").Value = 16.059769
$ws.Range("H21").Value = 32.119538
$ws.Range("I21").Value = 0.04831535056182164
$ws.Range("J21").Value = 0.032737474561192
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.06689266666666667
$ws.Range("N21").Value = 0.200678
$ws.Range("O21").Value = 0.001414591508805173
$ws.Range("P21").Value = 0.001478007612330618
$ws.Range("Q21").Value = 1.074280774460667
$ws.Range("R21").Value = 6.445684646764
$ws.Range("S21").Value = 0.00006834648464969813
$ws.Range("T21").Value = 0.00004838623660992174

# Row 22
$ws.Range("G22").Value = 0.5092873333333333
$ws.Range("H22").Value = 1.527862
$ws.Range("I22").Value = 0.001532176212914103
$ws.Range("J22").Value = 0.001557256002810873
$ws.Range("M22").Value = 32.21373866666666
$ws.Range("N22").Value = 96.641216
$ws.Range("O22").Value = 0.6812298485843321
$ws.Range("P22").Value = 0.7117693664123
$ws.Range("Q22").Value = 16.40604906224355
$ws.Range("R22").Value = 147.654441560192
$ws.Range("S22").Value = 0.00104376416952799
$ws.Range("T22").Value = 0.001108407118462446

# Row 23
$ws.Range("G23").Value = 0.5092873333333333
$ws.Range("H23").Value = 1.527862
$ws.Range("I23").Value = 0.001532176212914103
$ws.Range("J23").Value = 0.001557256002810873
$ws.Range("O23").Value = 0.1770389772624213
$ws.Range("P23").Value = 0.184975630381169
$ws.Range("Q23").Value = 4.263627251408
$ws.Range("R23").Value = 38.372645262672
$ws.Range("S23").Value = 0.0002712549097201226
$ws.Range("T23").Value = 0.0002880544107848007

# Row 24
$ws.Range("G24").Value = 0.5092873333333333
$ws.Range("H24").Value = 1.527862
$ws.Range("I24").Value = 0.001532176212914103
$ws.Range("J24").Value = 0.001557256002810873
$ws.Range("M24").Value = 0.5484013333333334
$ws.Range("N24").Value = 1.645204
$ws.Range("O24").Value = 0.01159714372603029
$ws.Range("P24").Value = 0.01211704340205096
$ws.Range("Q24").Value = 0.2792938526497778
$ws.Range("R24").Value = 2.513644673848
$ws.Range("S24").Value = 0.00001776886775476963
$ws.Range("T24").Value = 0.00001886933857416374

# Row 25
$ws.Range("G25").Value = 0.5092873333333333
$ws.Range("H25").Value = 1.527862
$ws.Range("I25").Value = 0.001532176212914103
$ws.Range("J25").Value = 0.001557256002810873
$ws.Range("M25").Value = 6.086836
$ws.Range("N25").Value = 12.173672
$ws.Range("O25").Value = 0.1287194389184112
$ws.Range("P25").Value = 0.08965995219214913
$ws.Range("Q25").Value = 3.099948474877333
$ws.Range("R25").Value = 18.599690849264
$ws.Range("S25").Value = 0.0001972208624504395
$ws.Range("T25").Value = 0.0001396234987629601

# Row 26
$ws.Range("G26").Value = 0.5092873333333333
$ws.Range("H26").Value = 1.527862
$ws.Range("I26").Value = 0.001532176212914103
$ws.Range("J26").Value = 0.001557256002810873
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.06689266666666667
$ws.Range("N26").Value = 0.200678
$ws.Range("O26").Value = 0.001414591508805173
$ws.Range("P26").Value = 0.001478007612330618
$ws.Range("Q26").Value = 0.03406758782622223
$ws.Range("R26").Value = 0.306608290436
$ws.Range("S26").Value = 0.000002167403460781557
$ws.Range("T26").Value = 0.000002849532370820838
